# Append new data rows (245-247) to the end of the sheet, matching the
# existing layout: column A = date (styled like the other date cells),
# B = nuovi positivi, C = somma mobile 7gg,
# D = somma mobile 7gg per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing date cell (A244) so the new
# date cells keep the same style (bold, centered, bordered, date format).
$formatSource = $ws.Range("A244")
$formatSource.Copy()

$newRows = @(
    @{ Row = 245; Date = 44319; B = 5; C = 27; D = 111.9310173285797 },
    @{ Row = 246; Date = 44320; B = 2; C = 23; D = 95.34864439101236 },
    @{ Row = 247; Date = 44321; B = 1; C = 23; D = 95.34864439101236 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.PasteSpecial(-4122)  # xlPasteFormats
    $cellA.Value = $r.Date

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}
